# Update "paises" (countries) workbook:
#  - Refresh the "datos actualizados" timestamp in A1
#  - Panama overtakes Oman in the ranking (rows 38/39 swap country + stats)
#  - Refresh numeric stats (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Muertes hoy, Muertes) for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 02:03"

# Columns B..H hold: Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes
$cols = @("B", "C", "D", "E", "F", "G", "H")

# row -> updated values for columns B..H
$updates = @{
    4   = @(5744688, 43757, 3085879, 2481470, 0, 1005, 177339)
    5   = @(3505097, 44684, 2653407, 739267, 0, 1234, 112423)
    9   = @(567059, 8639, 380730, 159295, 0, 200, 27034)
    16  = @(320884, 8225, 233651, 80716, 0, 187, 6517)
    74  = @(21045, 247, 15886, 4753, 0, 2, 406)
    114 = @(4229, 55, 3186, 962, 0, 1, 81)
    120 = @(3366, 71, 2383, 928, 0, 1, 55)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}

# Panama's numbers have grown past Oman's: Panama takes rank 42 (row 38),
# Oman drops to rank 43 (row 39). Oman's row keeps its previous stats.
$ws.Range("A38").Value = "Panama"
$ws.Range("B38").Value = 83855
$ws.Range("C38").Value = 101
$ws.Range("D38").Value = 59174
$ws.Range("E38").Value = 22837
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 17
$ws.Range("H38").Value = 1844

$ws.Range("A39").Value = "Oman"
$ws.Range("B39").Value = 83769
$ws.Range("C39").Value = 163
$ws.Range("D39").Value = 78386
$ws.Range("E39").Value = 4774
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 609

# Rows with a partial refresh (only a subset of columns changed)
$ws.Range("D135").Value = 1705
$ws.Range("E135").Value = 351

$ws.Range("B152").Value = 1297
$ws.Range("C152").Value = 12
$ws.Range("D152").Value = 1034
$ws.Range("E152").Value = 208

$ws.Range("B156").Value = 1169
$ws.Range("C156").Value = 2
$ws.Range("E156").Value = 18
